$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.235.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.894.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.67%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5065'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4040'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08301'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.112'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.406'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.888.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.328'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001100'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06469'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.26%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.927'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.261.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.187'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.100.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.278'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.121'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1044'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.014'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.717'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.91%  '
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.343'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2161'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.189'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.631'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6417'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.215'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5983'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.152'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.644'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.03%  '
